# 自动更新Excel文件 - 2026-01-09 23:13:36
# Decrement the "剩余" (remaining) counter in column E by 1 for every data
# row (rows 2-99), since one more day has elapsed. Row 36 is left untouched
# (its data was already inconsistent/skipped upstream). Row 95's counter
# would drop to 0, so it "rolls over": 剩余 resets to the 总天 (column D)
# value and 开始时间 (column F) advances by that many days (new cycle).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$skipRow = 36
$rolloverRow = 95

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq $skipRow) {
        continue
    }

    $remaining = $ws.Cells.Item($row, 5).Value2

    if ($row -eq $rolloverRow) {
        $totalDays = $ws.Cells.Item($row, 4).Value2
        $startDateNum = [int]$ws.Cells.Item($row, 6).Value2

        # F holds a plain YYYYMMDD integer (not a real Excel date serial),
        # so do real date math then re-encode back to YYYYMMDD.
        $y = [int]([math]::Floor($startDateNum / 10000))
        $m = [int]([math]::Floor(($startDateNum % 10000) / 100))
        $d = [int]($startDateNum % 100)
        $newDate = (Get-Date -Year $y -Month $m -Day $d).AddDays($totalDays)
        $newDateNum = ($newDate.Year * 10000) + ($newDate.Month * 100) + $newDate.Day

        $ws.Cells.Item($row, 5).Value2 = $totalDays
        $ws.Cells.Item($row, 6).Value2 = $newDateNum
    }
    else {
        $ws.Cells.Item($row, 5).Value2 = $remaining - 1
    }
}
